$d = $word.ActiveDocument

function Replace-Exact($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# --- Table 1 (Environmental characteristics): Chlorophyll-a column values ---
# Replace the most specific (longer / prefixed) strings first so that
# "-0.1 (-1.2, 1.1)" is not accidentally matched by "0.1 (-1.2, 1.1)".
Replace-Exact "-0.1 (-1.2, 1.1)" "1.9 (0.1, 11.7)"
Replace-Exact "-0.1 (-1.5, 1.7)" "4.1 (0, 47.6)"
Replace-Exact "0 (-0.7, 0.9)" "1.5 (0.2, 7.3)"
Replace-Exact "0.3 (-1.7, 1.3)" "3.8 (0, 20.8)"
Replace-Exact "0.1 (-1.2, 1.1)" "2.6 (0.1, 14)"
Replace-Exact "0.2 (-1.1, 1.6)" "4.9 (0.1, 43.9)"
Replace-Exact "0.1 (-1.3, 2.1)" "6 (0, 116.5)"

# --- Table 2 caption (#tab:araminmod): expand wording per SRA comments ---
Replace-Exact "relationships grouped by month (Figure 5, bottom row). Separate models were also run with and without 2016 data because of missing April observations in 2016." "relationships grouped by month (Figure 5, bottom row). Models in the first two columns included a year variable as a fixed effect and models in the second two columns included a month variable as a fixed effect, in addition to minimum aragonite saturation state as a fixed effect for all models. Separate models were also run with and without 2016 cohort-year data because of missing April observations in 2017 (the final month for the 2016 cohort)."
